# Updated cryptos list - price (D) and volume/%-change (E) refresh,
# plus an Aave/Maker row-order swap (rows 47-48).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as literal text (e.g. "309.78", "42.991.23" using
# "." as a thousands separator) in the source data. Force Text format before
# assigning so COM does not silently reinterpret these as locale numbers
# (which would corrupt values like "1.00" -> 1 or introduce float drift).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.000.40"
$ws.Range("E2").Value = "  +0.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.304.10"
$ws.Range("E3").Value = "  +0.05%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.78"
$ws.Range("E5").Value = "  -2.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.12"
$ws.Range("E6").Value = "  +1.23%  "

$ws.Range("E7").Value = "  -0.46%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.72"
$ws.Range("E10").Value = "  -0.29%  "

$ws.Range("E11").Value = "  +0.74%  "

$ws.Range("E12").Value = "  -2.70%  "

$ws.Range("E13").Value = "  +0.02%  "

$ws.Range("E14").Value = "  -1.07%  "

$ws.Range("E15").Value = "  -0.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.652.42"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.304.94"
$ws.Range("E17").Value = "  -0.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.812.95"
$ws.Range("E18").Value = "  +0.44%  "

$ws.Range("E19").Value = "  -3.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.81"
$ws.Range("E20").Value = "  -0.85%  "

$ws.Range("E21").Value = "  -1.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.49"
$ws.Range("E22").Value = "  -0.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.46"
$ws.Range("E23").Value = "  -2.57%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.52"
$ws.Range("E24").Value = "  +0.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.24"
$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.71"
$ws.Range("E27").Value = "  +16.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.94"
$ws.Range("E28").Value = "  +0.21%  "

$ws.Range("E29").Value = "  +1.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.75"
$ws.Range("E30").Value = "  +0.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.21"
$ws.Range("E31").Value = "  -1.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.85"
$ws.Range("E32").Value = "  +0.06%  "

$ws.Range("E33").Value = "  -2.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.83"
$ws.Range("E34").Value = "  +6.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("E35").Value = "  -0.88%  "

$ws.Range("E36").Value = "  -1.02%  "

$ws.Range("E37").Value = "  +1.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0358"
$ws.Range("E38").Value = "  +0.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.78"
$ws.Range("E39").Value = "  +1.14%  "

$ws.Range("E40").Value = "  -3.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "108.51"
$ws.Range("E41").Value = "  +12.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.56"
$ws.Range("E42").Value = "  -3.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.27"
$ws.Range("E43").Value = "  +0.90%  "

$ws.Range("E44").Value = "  +0.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("E46").Value = "  -1.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "75.99"
$ws.Range("E49").Value = "  -5.40%  "

$ws.Range("E50").Value = "  -0.44%  "

$ws.Range("E51").Value = "  -2.38%  "

# Rows 47 and 48 swapped rank order: Maker now above Aave.
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.700.79"
$ws.Range("E47").Value = "  +2.92%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "111.61"
$ws.Range("E48").Value = "  -5.14%  "

